# Update the font-size preset descriptions on the "Font Size Presets" slide
# (slide 3) to include the point value inline and a short usage hint,
# replacing the old "LABEL: N (Npt)" wording with "LABEL: Npt - For ...".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$content = $s.Shapes.Item(2)
$tf = $content.TextFrame.TextRange

$tf.Paragraphs(1).Runs(1).Text = "TITLE: 44pt - For main titles"
$tf.Paragraphs(2).Runs(1).Text = "SUBTITLE: 32pt - For subtitles"
$tf.Paragraphs(3).Runs(1).Text = "HEADING: 28pt - For section headers"
$tf.Paragraphs(4).Runs(1).Text = "BODY: 18pt - For regular content"
$tf.Paragraphs(5).Runs(1).Text = "SMALL: 14pt - For smaller text"
$tf.Paragraphs(6).Runs(1).Text = "CAPTION: 12pt - For captions"
